$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.933.71'
$ws.Range("E2").Value = '  -1.10%  '

# Row 3
$ws.Range("D3").Value = '3.503.62'
$ws.Range("E3").Value = '  +0.26%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.31%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.39%  '

# Row 7
$ws.Range("E7").Value = '  +4.67%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("E9").Value = '  +0.10%  '

# Row 10
$ws.Range("E10").Value = '  +3.82%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '55.36'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.48%  '

# Row 12
$ws.Range("E12").Value = '  +1.81%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.54%  '

# Row 14
$ws.Range("D14").Value = '4.068.14'
$ws.Range("E14").Value = '  +0.21%  '

# Row 15
$ws.Range("D15").Value = '3.502.95'
$ws.Range("E15").Value = '  +0.20%  '

# Row 16
$ws.Range("E16").Value = '  +0.07%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.82%  '

# Row 18
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.48%  '

# Row 19
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '65.907.62'
$ws.Range("E19").Value = '  -1.04%  '

# Row 20
$ws.Range("E20").Value = '  +1.24%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '414.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.34%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.36%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.03%  '

# Row 24
$ws.Range("E24").Value = '  +0.88%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.28'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +12.17%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.91%  '

# Row 27
$ws.Range("E27").Value = '  -2.43%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.03%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '30.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.75%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '624.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.60%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.50'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.81%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.66'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.27%  '

# Row 33
$ws.Range("E33").Value = '  -0.56%  '

# Row 34
$ws.Range("E34").Value = '  +14.27%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '59.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.22%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.20%  '

# Row 37
$ws.Range("E37").Value = '  -2.05%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.24'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.63%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.97%  '

# Row 40
$ws.Range("D40").Value = '3.289.91'
$ws.Range("E40").Value = '  +10.01%  '

# Row 41
$ws.Range("E41").Value = '  -2.97%  '

# Row 42
$ws.Range("E42").Value = '  +0.01%  '

# Row 43
$ws.Range("E43").Value = '  +0.04%  '

# Row 44
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.16%  '

# Row 45
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0417'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.42%  '

# Row 46
$ws.Range("E46").Value = '  -4.82%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.11%  '

# Row 48
$ws.Range("E48").Value = '  +1.83%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.11'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.83%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.57'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.74%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.09%  '
